$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert a new blank column before column D (shifts D:K -> E:L)
$ws.Range("D:D").Insert()

# Step 2: Copy number formats from column E into new column D so the new
#         period column matches the look of the existing data (dates, numbers).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Step 3: Populate the new column D with the new period values.
# Numeric values (including zeros for blank line items):
$newDValues = @{
  7 = 43465
  8 = 1571300
  9 = 642200
  10 = 929100
  13 = 0
  14 = 9100
  15 = 489800
  17 = 1175900
  18 = 395500
  20 = 9900
  21 = 896400
  22 = 173600
  23 = 231800
  24 = 2600
  25 = 0
  26 = 229200
  27 = 218900
  28 = 0
  30 = 0
  31 = 0
  32 = -9900
  33 = 218900
  34 = 0
  35 = 218900
  38 = 43465
  41 = 34300
  42 = 0
  43 = 0
  44 = 0
  45 = 0
  46 = 0
  47 = 44200
  48 = 11107500
  50 = 0
  51 = 0
  52 = 17400
  53 = 0
  54 = 11323800
  58 = 0
  59 = 413900
  60 = 0
  61 = 4528300
  62 = 0
  63 = 0
  64 = 0
  65 = 0
  66 = 5164500
  68 = 0
  69 = 0
  70 = 9400
  71 = 0
  72 = -989300
  73 = 0
  74 = 0
  75 = 0
  76 = 6149900
  77 = 0
  80 = 43465
  81 = 218900
  83 = 491000
  84 = 0
  85 = 0
  86 = 0
  87 = 0
  88 = 0
  91 = -129500
  92 = 0
  93 = 0
  96 = -423500
  97 = 0
  98 = 0
  99 = 0
  101 = 0
}
foreach ($row in $newDValues.Keys) {
  $ws.Cells.Item($row, 4).Value2 = $newDValues[$row]
}

# Rows whose new-period value is the text "NA":
$newDTextRows = @(12, 29, 49, 57)
foreach ($row in $newDTextRows) {
  $ws.Cells.Item($row, 4).Value2 = "NA"
}

# Step 4: A handful of cash-flow subtotal rows were recalculated along with
#         the new period, so their D (new) and shifted E/F values differ from
#         a pure shift of the prior data. Apply the corrected values explicitly.
$specialRows = @{
  89 = @{ D = 734300; E = 660800; F = 485000 }
  94 = @{ D = -366400; E = -294200; F = -649100 }
  100 = @{ D = -405100; E = -399500; F = 222300 }
  102 = @{ D = -37200; E = -32900; F = 58200 }
}
foreach ($row in $specialRows.Keys) {
  $ws.Cells.Item([int]$row, 4).Value2 = $specialRows[$row].D
  $ws.Cells.Item([int]$row, 5).Value2 = $specialRows[$row].E
  $ws.Cells.Item([int]$row, 6).Value2 = $specialRows[$row].F
}

Write-Host "Edit complete"
